# Update gh-pages output data (F/G columns) for sheets "展览" and "全部类型"
$wb = $excel.ActiveWorkbook

# Row -> new F (want-to-go count) value; only row 2 also changes G (min ticket price)
$updates = @(
    @{ Row = 2;  F = 144;  G = 60 },
    @{ Row = 7;  F = 1318 },
    @{ Row = 8;  F = 1559 },
    @{ Row = 10; F = 428 },
    @{ Row = 13; F = 165 },
    @{ Row = 16; F = 278 },
    @{ Row = 18; F = 331 },
    @{ Row = 19; F = 1761 },
    @{ Row = 23; F = 686 },
    @{ Row = 26; F = 4239 },
    @{ Row = 29; F = 1113 },
    @{ Row = 30; F = 495 },
    @{ Row = 32; F = 610 },
    @{ Row = 34; F = 309 },
    @{ Row = 36; F = 153 }
)

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates) {
        $ws.Cells.Item($u.Row, 6).Value = $u.F
        if ($u.ContainsKey("G")) {
            $ws.Cells.Item($u.Row, 7).Value = $u.G
        }
    }
}
